$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C ("Förändrad") date value from 45720 (2025-03-04) to
# 45721 (2025-03-05) for every data row (rows 2 through 42).
for ($row = 2; $row -le 42; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45720) {
        $cell.Value2 = 45721
    }
}
